$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K")

# Insert a new blank row at position 7; this pushes the existing
# "Giorni settimana" (old row 7) and "Mesi anno" (old row 8) rows down
# to rows 8 and 9 respectively.
$ws.Rows("7:7").Insert()

# --- New row 10: "Date:" example (sequential dates with date formatting) ---
$ws.Range("A10").Value = "Date:"
$dateValues = @(43160, 43161, 43162, 43163, 43164, 43165, 43166, 43167, 43168, 43169)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "10").Value = $dateValues[$i]
}
$ws.Range("B10").NumberFormat = "mm-dd-yy"
$ws.Range("B10").Copy()
$ws.Range("C10:K10").PasteSpecial(-4122)

# --- New row 11: "Orari:" example (half-hour increments with time formatting) ---
$ws.Range("A11").Value = "Orari:"
$timeValues = @(0.22916666666666666, 0.25, 0.27083333333333298, 0.29166666666666702, 0.3125, 0.33333333333333398, 0.35416666666666702, 0.375, 0.39583333333333398, 0.41666666666666702)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "11").Value = $timeValues[$i]
}
$ws.Range("B11").NumberFormat = "h:mm"
$ws.Range("B11").Copy()
$ws.Range("C11:K11").PasteSpecial(-4122)

# --- New row 7: "Incremento decimali:" example (0.5 step increments) ---
$ws.Range("A7").Value = "Incremento decimali:"
$decimalValues = @(1.5, 2.5, 3.5, 4.5, 5.5, 6.5, 7.5, 8.5, 9.5, 10.5)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "7").Value = $decimalValues[$i]
}

$excel.CutCopyMode = $false

# Update the selection to match the authored state (B11:K11, active cell B11)
[void]$ws.Range("B11:K11").Select()

Write-Output "done"
